# "Generate Report for Handoff"
#
# For the six rows whose source file already has a zh-cn handoff xliff
# generated (1953d55a, 2e73bd87, 31c729d1, 633a229a, 6818f53b, 6c64d27c —
# i.e. table rows 8,9,10,12,13,14 on each localized-language sheet), mark
# the row's Priority as "ht" (handoff type) and refresh the handoff
# generation timestamps to reflect the new report run.

$wb = $excel.ActiveWorkbook

$rows = @(8, 9, 10, 12, 13, 14)

# Overview sheet: column G = "Latest HO Xliff Generate Date" (zh-cn value)
$wsOverview = $wb.Worksheets.Item("Overview")
foreach ($r in $rows) {
    $wsOverview.Range("G$r").Value = "2016-09-05 06:23:41"
}

# zh-cn sheet: column E = "Priority", column H = "Latest Handoff Datetime"
$wsZhCn = $wb.Worksheets.Item("zh-cn")
foreach ($r in $rows) {
    $wsZhCn.Range("E$r").Value = "ht"
    $wsZhCn.Range("H$r").Value = "2016-09-05 06:23:35"
}

# de-de sheet: column E = "Priority", column H = "Latest Handoff Datetime"
$wsDeDe = $wb.Worksheets.Item("de-de")
foreach ($r in $rows) {
    $wsDeDe.Range("E$r").Value = "ht"
    $wsDeDe.Range("H$r").Value = "2016-09-05 06:23:41"
}
